$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column F ("MS_DEF") carrying the MS definition list for each
# mapping row. Copy the header style from the existing E1 header cell so
# the new header matches the other headers (bold, centered, bordered).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "MS_DEF"

$ws.Range("F2").Value = '["A physical quality that inheres in a bearer by virtue of the bearer''s amount of force per unit area it exerts."]'
$ws.Range("F3").Value = '[''A physical quality of the thermal energy of a system. [PATO]'']'
$ws.Range("F4").Value = '["A physical quality which inheres in a bearer by virtue of some influence is exerted by the bearer''s mass per unit size."]'
$ws.Range("F5").Value = "[]"
$ws.Range("F6").Value = "[]"
$ws.Range("F7").Value = "[]"
$ws.Range("F8").Value = "[]"
$ws.Range("F9").Value = "[]"
$ws.Range("F10").Value = "[]"
$ws.Range("F11").Value = "[]"
$ws.Range("F12").Value = "[]"
$ws.Range("F13").Value = "[]"
$ws.Range("F14").Value = "[]"
$ws.Range("F15").Value = "[]"
$ws.Range("F16").Value = "[]"
$ws.Range("F17").Value = "[]"
$ws.Range("F18").Value = "[]"
$ws.Range("F19").Value = "[]"
$ws.Range("F20").Value = "[]"
$ws.Range("F21").Value = '["A physical quality that inheres in a bearer by virtue of the proportion of the bearer''s amount of matter. [PATO]"]'
$ws.Range("F22").Value = "[]"
$ws.Range("F23").Value = "[]"
$ws.Range("F24").Value = "[]"
$ws.Range("F25").Value = "[]"
$ws.Range("F26").Value = "[]"
$ws.Range("F27").Value = "[]"
